$d = $word.ActiveDocument

# --- Edit 1: "Full Bayesian methods..." paragraph -----------------------------------------
# Removes the two proofErr (gramStart/gramEnd) markers around "second, and" and merges the
# surrounding runs into a single run (text itself is unchanged).
$d.Content.Find.Execute(", such as analyzing a scientific study. However, production systems often need to perform estimation in a fraction of a second, and run them thousands or millions of times each day.", $true, $false, $false, $false, $false, $true, 1, $false, ", such as analyzing a scientific study. However, production systems often need to perform estimation in a fraction of a second, and run them thousands or millions of times each day.", 2) | Out-Null

# --- Edit 2: "Empirical Bayes estimation, where a beta distribution fiton ..." paragraph ----
# Fixes the "fiton" -> "fit on" typo, removes the gramStart/gramEnd markers around
# "As long as" and merges the trailing runs together.
$d.Content.Find.Execute("fiton all observations is then used to improve each individually. As long as you have a lot of examples, you don’t need to bring in prior expectations.", $true, $false, $false, $false, $false, $true, 1, $false, "fit on all observations is then used to improve each individually. As long as you have a lot of examples, you don’t need to bring in prior expectations.", 2) | Out-Null

# --- Edit 3: "Imagine the player gets a single hit..." paragraph --------------------------
# Removes the gramStart/gramEnd markers around "have to" and merges the runs into one.
$d.Content.Find.Execute("Imagine the player gets a single hit. His record for the season is now “1 hit; 1 at bat.” We have to then update our probabilities- we want to shift this entire curve over just a bit to reflect our new information. According to that our new parameters will be:", $true, $false, $false, $false, $false, $true, 1, $false, "Imagine the player gets a single hit. His record for the season is now “1 hit; 1 at bat.” We have to then update our probabilities- we want to shift this entire curve over just a bit to reflect our new information. According to that our new parameters will be:", 2) | Out-Null

# --- Edit 4: mark the two chart pictures as NoProof (w:rPr/w:noProof) ----------------------
$d.InlineShapes.Item(25).Range.NoProofing = 1
$d.InlineShapes.Item(26).Range.NoProofing = 1

# --- Edit 5: "Notice that crossover point..." paragraph -----------------------------------
# Merges the five runs that make up the sentence into a single run.
$d.Content.Find.Execute("Notice that crossover point: to have a PEP less than 50%, you need to have a shrunken batting average greater than .300.", $true, $false, $false, $false, $false, $true, 1, $false, "Notice that crossover point: to have a PEP less than 50%, you need to have a shrunken batting average greater than .300.", 2) | Out-Null

Write-Output "edits applied"
